# Increase font sizes throughout the resume.
#
# Every run-level font size (w:sz) in the document grows by 1pt (2
# half-points), except the name header ("Michael L. Welles"), which is
# bolded at 16pt and grows by 2pt (to 18pt) -- a special case confirmed
# by the diff (w:sz 32 -> 36, a jump of 4 half-points instead of 2).
#
# Each paragraph in this resume uses a single, uniform font size across
# all of its runs, so adjusting the font size paragraph-by-paragraph
# reproduces the diff exactly without disturbing any other run property
# (bold/italic/color/etc. are left untouched).
#
# Note: we rebuild each paragraph's range via $d.Range(start, end)
# rather than using $p.Range directly -- assigning Font.Size on the
# Paragraph's own Range object also stamps the paragraph-mark's run
# properties (w:pPr/w:rPr/w:sz), which is not part of the target diff.
# A Range constructed from the same start/end offsets edits only the
# visible run content, matching the original document structure.

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $pr = $p.Range
    $r = $d.Range($pr.Start, $pr.End)
    $current = $r.Font.Size

    if ($current -eq 16) {
        $r.Font.Size = 18
    } else {
        $r.Font.Size = $current + 1
    }
}
